$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Pandas ,GeoPandas , sqlite3, psycopg2"
$ws.Range("B3").Value = "numpy, datatime, random"

$ws.Range("B3").Select()
